$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlLeft = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignLeft
$xlRight = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignRight

# --- Row 1 header: column C label changes from "Co" to "Thuoc" (style unchanged: centered) ---
$ws.Cells.Item(1, 3).Value2 = 'Thuộc'

# --- Data rows: update/insert A, B, C values ---
$ws.Cells.Item(2, 1).Value2 = 'PQLCL'
$ws.Cells.Item(2, 2).Value2 = 'Phòng quản lý chất lượng'
$ws.Cells.Item(2, 3).Value2 = 0

$ws.Cells.Item(3, 1).Value2 = 'MT'
$ws.Cells.Item(3, 2).Value2 = 'Khoa Khoa học và Kỹ thuật máy tính'
$ws.Cells.Item(3, 3).Value2 = 0

$ws.Cells.Item(4, 1).Value2 = 'MT-HTTT'
$ws.Cells.Item(4, 2).Value2 = 'Hệ thống thông tin'
$ws.Cells.Item(4, 3).Value2 = 'MT'

$ws.Cells.Item(5, 1).Value2 = 'MT-KHMT'
$ws.Cells.Item(5, 2).Value2 = 'Khoa học máy tính'
$ws.Cells.Item(5, 3).Value2 = 'MT'

$ws.Cells.Item(6, 1).Value2 = 'MT-KTMT'
$ws.Cells.Item(6, 2).Value2 = 'Kỹ thuật máy tính'
$ws.Cells.Item(6, 3).Value2 = 'MT'

$ws.Cells.Item(7, 1).Value2 = 'MT-HTM'
$ws.Cells.Item(7, 2).Value2 = 'Hệ thống mạng'
$ws.Cells.Item(7, 3).Value2 = 'MT'

$ws.Cells.Item(8, 1).Value2 = 'MT-KPTT'
$ws.Cells.Item(8, 2).Value2 = 'Khám phá tri thức'
$ws.Cells.Item(8, 3).Value2 = 'MT'

$ws.Cells.Item(9, 1).Value2 = 'KHUD'
$ws.Cells.Item(9, 2).Value2 = 'Khoa Khoa học ứng dụng'
$ws.Cells.Item(9, 3).Value2 = 0

$ws.Cells.Item(10, 1).Value2 = 'KHUD-TUD'
$ws.Cells.Item(10, 2).Value2 = 'Toán ứng dụng'
$ws.Cells.Item(10, 3).Value2 = 'KHUD'

$ws.Cells.Item(11, 1).Value2 = 'KHUD-CKT'
$ws.Cells.Item(11, 2).Value2 = 'Cơ Kỹ Thuật'
$ws.Cells.Item(11, 3).Value2 = 'KHUD'

$ws.Cells.Item(12, 1).Value2 = 'KHUD-CHUD'
$ws.Cells.Item(12, 2).Value2 = 'Cơ học ứng dụng'
$ws.Cells.Item(12, 3).Value2 = 'KHUD'

$ws.Cells.Item(13, 1).Value2 = 'KHUD-VLKT'
$ws.Cells.Item(13, 2).Value2 = 'Vật lý kỹ thuật'
$ws.Cells.Item(13, 3).Value2 = 'KHUD'

$ws.Cells.Item(14, 1).Value2 = 'KHUD-VLUD'
$ws.Cells.Item(14, 2).Value2 = 'Vật lý ứng dụng'
$ws.Cells.Item(14, 3).Value2 = 'KHUD'

$ws.Cells.Item(15, 1).Value2 = 'KHUD-VLDC'
$ws.Cells.Item(15, 2).Value2 = 'Vật lý đại cương'
$ws.Cells.Item(15, 3).Value2 = 'KHUD'

$ws.Cells.Item(16, 1).Value2 = 'KHUD-LLCT'
$ws.Cells.Item(16, 2).Value2 = 'Lý luận chính trị'
$ws.Cells.Item(16, 3).Value2 = 'KHUD'

# --- Formatting (applied in the same order the styles first appear, so the
#     resulting style indices line up with the target workbook) ---
# C2: right-aligned (no wrap)
$ws.Range("C2").HorizontalAlignment = $xlRight

# Rows 3-16 column C: right-aligned + wrap text
$cRange = $ws.Range("C3:C16")
$cRange.HorizontalAlignment = $xlRight
$cRange.WrapText = $true

# Row 2 (PQLCL): A2/B2 left-aligned
$ws.Range("A2:B2").HorizontalAlignment = $xlLeft

# --- sheet view: update selection to C2 ---
[void]$ws.Range("C2").Select()
